# Updated cryptos list with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures scraped
# from coinranking.com for each coin row in the sheet.
#
# Note: several Price values are plain numeric-looking strings (e.g. "213.13").
# Excel's Range.Value setter auto-coerces those into numbers, which would
# silently flip the cell from text to a numeric type. A leading apostrophe
# (the classic "quote prefix") forces Excel to keep them as text, matching
# the original inline-string cell type, exactly as typing '213.13 into the
# cell in the UI would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r -> @{ D = <new price or $null>; E = <new volume text> }
$updates = [ordered]@{
    2  = @{ D = "26.149.04";  E = "  +3.49%  " }
    3  = @{ D = "1.604.06";   E = "  +3.42%  " }
    4  = @{ D = $null;        E = "  -0.22%  " }
    5  = @{ D = "'213.13";    E = "  +3.11%  " }
    6  = @{ D = $null;        E = "  -0.21%  " }
    7  = @{ D = $null;        E = "  +1.92%  " }
    8  = @{ D = $null;        E = "  +2.67%  " }
    9  = @{ D = $null;        E = "  +1.30%  " }
    10 = @{ D = "'18.08";     E = "  +2.42%  " }
    11 = @{ D = "'0.0818";    E = "  +5.01%  " }
    12 = @{ D = "1.828.53";   E = "  +3.47%  " }
    13 = @{ D = "1.604.93";   E = "  +3.51%  " }
    14 = @{ D = $null;        E = "  +0.88%  " }
    15 = @{ D = "'0.512";     E = "  +1.74%  " }
    16 = @{ D = "26.163.13";  E = $null }
    17 = @{ D = "'60.61";     E = "  +3.59%  " }
    18 = @{ D = "0.0₃0723";   E = "  +2.68%  " }
    19 = @{ D = $null;        E = "  -0.22%  " }
    20 = @{ D = "'204.34";    E = "  +10.13%  " }
    21 = @{ D = $null;        E = "  +3.47%  " }
    22 = @{ D = $null;        E = "  +0.78%  " }
    23 = @{ D = $null;        E = "  +2.98%  " }
    24 = @{ D = "'1.85";      E = "  +13.84%  " }
    25 = @{ D = "'141.95";    E = "  +2.34%  " }
    26 = @{ D = $null;        E = "  -0.27%  " }
    27 = @{ D = "'0.123";     E = "  -4.48%  " }
    28 = @{ D = $null;        E = "  +2.67%  " }
    29 = @{ D = "'6.43";      E = "  +0.87%  " }
    30 = @{ D = $null;        E = "  +1.74%  " }
    31 = @{ D = $null;        E = "  +2.31%  " }
    32 = @{ D = "'3.11";      E = "  +3.18%  " }
    33 = @{ D = $null;        E = "  +1.26%  " }
    34 = @{ D = $null;        E = "  +2.03%  " }
    35 = @{ D = $null;        E = "  +1.58%  " }
    36 = @{ D = $null;        E = "  +11.64%  " }
    37 = @{ D = "1.125.69";   E = "  +4.05%  " }
    38 = @{ D = $null;        E = "  +0.09%  " }
    39 = @{ D = $null;        E = "  +3.10%  " }
    40 = @{ D = $null;        E = "  +2.70%  " }
    41 = @{ D = $null;        E = "  +0.29%  " }
    42 = @{ D = "'0.781";     E = "  -1.89%  " }
    43 = @{ D = "'5.16";      E = "  +2.40%  " }
    44 = @{ D = "1.740.73";   E = "  +3.48%  " }
    45 = @{ D = "'92.92";     E = "  +0.24%  " }
    46 = @{ D = $null;        E = "  +4.70%  " }
    48 = @{ D = $null;        E = "  +0.53%  " }
    49 = @{ D = "'0.409";     E = "  +1.15%  " }
    50 = @{ D = $null;        E = "  -0.20%  " }
    51 = @{ D = "0.0₇0927";   E = "  -16.66%  " }
}

foreach ($r in $updates.Keys) {
    $row = $updates[$r]
    if ($null -ne $row.D) {
        $ws.Range("D$r").Value = $row.D
    }
    if ($null -ne $row.E) {
        $ws.Range("E$r").Value = $row.E
    }
}
